$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-07 Wednesday" "2024-02-08 Thursday"

Replace-Text "971×6=5826" "928×2=1856"
Replace-Text "936×6=5616" "640×2=1280"
Replace-Text "869×6=5214" "120×8=960"
Replace-Text "905×5=4525" "489×9=4401"
Replace-Text "725×7=5075" "722×5=3610"

Replace-Text "901×6=5406" "961×5=4805"
Replace-Text "790×9=7110" "927×2=1854"
Replace-Text "784×3=2352" "123×6=738"
Replace-Text "870×5=4350" "221×3=663"
Replace-Text "744×2=1488" "207×4=828"

Replace-Text "756×9=6804" "226×3=678"
Replace-Text "267×4=1068" "807×4=3228"
Replace-Text "456×5=2280" "798×4=3192"
Replace-Text "639×3=1917" "440×6=2640"
Replace-Text "633×4=2532" "623×6=3738"

Replace-Text "933×7=6531" "357×8=2856"
Replace-Text "754×2=1508" "315×9=2835"
Replace-Text "487×5=2435" "973×4=3892"
Replace-Text "810×4=3240" "162×6=972"
Replace-Text "419×8=3352" "813×9=7317"

Replace-Text "164×3=492" "852×3=2556"
Replace-Text "766×5=3830" "945×9=8505"
Replace-Text "839×8=6712" "905×8=7240"
Replace-Text "918×3=2754" "750×2=1500"
Replace-Text "602×6=3612" "981×8=7848"
